$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "248.07"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "4"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.75"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "4"

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.483"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "4"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05639"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "4"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.382"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "4"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.448"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "4"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8011"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "4"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.033"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "4"

# Row 10
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "One"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01152"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "4"

# Row 11
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1427"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "4"

# Row 12
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07247"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "4"

# Row 13
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03112"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "4"

# Row 14
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.02940"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "4"

# Row 15
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09284"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "4"

# Row 16
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001650"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "4"

# Row 17
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "MCDex"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.257"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "16MCDexMCB"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "4"

# Row 18
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04732"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "4"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006374"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "4"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005025"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "4"

# Row 21
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "4"

# Row 22
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "4"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0003206"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "4"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.167"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "4"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.086"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "4"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.3271"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "4"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1309"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "4"

# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "4"

# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "4"

# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "4"

# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "4"

# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "4"

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "4"

# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "4"

# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "4"

# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "4"

# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "4"

# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "4"

# Row 39
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "4"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04075"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "4"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006914"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "4"

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003507"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "4"

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1040"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "4"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009123"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "4"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005646"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "4"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "4"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7867"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "4"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.01667"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "4"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002104"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "4"

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.01012"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "4"

# Row 51
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "4"
